$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3999.96
$ws.Range("J40").Value = 3999.96
$ws.Range("L40").Value = 3999.96
$ws.Range("N40").Value = -4349.96

$ws.Range("H103").Value = 582.94446
$ws.Range("J103").Value = 566.9286
$ws.Range("L103").Value = 1700.7858
$ws.Range("N103").Value = -2872.7858

$ws.Range("H111").Value = 3860.5557
$ws.Range("I111").Value = 679
$ws.Range("J111").Value = 6405.8
$ws.Range("K111").Value = 2037
$ws.Range("L111").Value = 19217.4
$ws.Range("M111").Value = 1030
$ws.Range("N111").Value = -25351.4

$ws.Range("H116").Value = 14414.538
$ws.Range("I116").Value = 20749.125
$ws.Range("J116").Value = 4279.2
$ws.Range("K116").Value = 20749.125
$ws.Range("L116").Value = 4279.2
$ws.Range("M116").Value = -17307.125
$ws.Range("N116").Value = -11163.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3791.724
$ws.Range("I61").Value = 3479.111
$ws.Range("K61").Value = 3479.111
$ws.Range("M61").Value = -3267.111

$ws.Range("H132").Value = 3164.56
$ws.Range("I132").Value = 3141.5454
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 9424.636200000001
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -6894.636200000001
$ws.Range("N132").Value = -15059.9999

$ws.Range("H136").Value = 3791.724
$ws.Range("I136").Value = 3479.111
$ws.Range("K136").Value = 10437.333
$ws.Range("M136").Value = -7887.332999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2297.762
$ws.Range("I99").Value = 1975.1875
$ws.Range("K99").Value = 1975.1875
$ws.Range("M99").Value = -477.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 3636666.2
$ws.Range("I23").Value = 5452499.5
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 5452499.5
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = -5452259.5
$ws.Range("N23").Value = -5480

$ws.Range("H27").Value = 3636666.2
$ws.Range("I27").Value = 5452499.5
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 5452499.5
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -5452307.5
$ws.Range("N27").Value = -5384

$ws.Range("H58").Value = 2649.1694
$ws.Range("I58").Value = 2441.2642
$ws.Range("K58").Value = 2441.2642
$ws.Range("M58").Value = -2238.2642

$ws.Range("H105").Value = 3780
$ws.Range("I105").Value = 3400
$ws.Range("J105").Value = 4350
$ws.Range("K105").Value = 3400
$ws.Range("L105").Value = 4350
$ws.Range("M105").Value = -1653
$ws.Range("N105").Value = -7844

$ws.Range("H122").Value = 3506.724
$ws.Range("I122").Value = 2711.7222
$ws.Range("J122").Value = 4807.636
$ws.Range("K122").Value = 8135.1666
$ws.Range("L122").Value = 14422.908
$ws.Range("M122").Value = -5685.1666
$ws.Range("N122").Value = -19322.908

$ws.Range("H131").Value = 3499.5
$ws.Range("J131").Value = 3499.5
$ws.Range("L131").Value = 3499.5
$ws.Range("N131").Value = -13579.5

$ws.Range("H132").Value = 4173.1816
$ws.Range("I132").Value = 3403.0344
$ws.Range("K132").Value = 10209.1032
$ws.Range("M132").Value = -7679.1032

$ws.Range("H136").Value = 2649.1694
$ws.Range("I136").Value = 2441.2642
$ws.Range("K136").Value = 7323.792600000001
$ws.Range("M136").Value = -4773.792600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9681.333000000001
$ws.Range("I3").Value = 9681.333000000001
$ws.Range("K3").Value = 29043.999
$ws.Range("M3").Value = -28931.999

$ws.Range("H133").Value = 4409.846
$ws.Range("I133").Value = 4041
$ws.Range("K133").Value = 12123
$ws.Range("M133").Value = -7063

$ws.Range("H138").Value = 42601096
$ws.Range("I138").Value = 1823.3334
$ws.Range("K138").Value = 5470.0002
$ws.Range("M138").Value = -330.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3643.75
$ws.Range("I80").Value = 3487.5
$ws.Range("K80").Value = 3487.5
$ws.Range("M80").Value = -2489.5

$ws.Range("H83").Value = 3643.75
$ws.Range("I83").Value = 3487.5
$ws.Range("K83").Value = 17437.5
$ws.Range("M83").Value = -12445.5

$ws.Range("H97").Value = 695.381
$ws.Range("I97").Value = 530.15
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 530.15
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -34.14999999999998
$ws.Range("N97").Value = -4992

$ws.Range("H113").Value = 83353416
$ws.Range("I113").Value = 142873570
$ws.Range("K113").Value = 142873570
$ws.Range("M113").Value = -142871400

$ws.Range("H132").Value = 3047
$ws.Range("I132").Value = 2565
$ws.Range("K132").Value = 7695
$ws.Range("M132").Value = -5165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7494.4614
$ws.Range("I7").Value = 7533.3
$ws.Range("J7").Value = 7365
$ws.Range("K7").Value = 7533.3
$ws.Range("L7").Value = 7365
$ws.Range("M7").Value = -7421.3
$ws.Range("N7").Value = -7589

$ws.Range("H16").Value = 1710.3334
$ws.Range("I16").Value = 1710.3334
$ws.Range("K16").Value = 1710.3334
$ws.Range("M16").Value = -1540.3334

$ws.Range("H68").Value = 2039.4546
$ws.Range("I68").Value = 1890.1428
$ws.Range("J68").Value = 2300.75
$ws.Range("K68").Value = 1890.1428
$ws.Range("L68").Value = 2300.75
$ws.Range("M68").Value = -1141.1428
$ws.Range("N68").Value = -3798.75

$ws.Range("H71").Value = 2039.4546
$ws.Range("I71").Value = 1890.1428
$ws.Range("J71").Value = 2300.75
$ws.Range("K71").Value = 9450.714
$ws.Range("L71").Value = 11503.75
$ws.Range("M71").Value = -5706.714
$ws.Range("N71").Value = -18991.75

$ws.Range("H126").Value = 7494.4614
$ws.Range("I126").Value = 7533.3
$ws.Range("J126").Value = 7365
$ws.Range("K126").Value = 22599.9
$ws.Range("L126").Value = 22095
$ws.Range("M126").Value = -20129.9
$ws.Range("N126").Value = -27035

$ws.Range("H136").Value = 5930.9546
$ws.Range("I136").Value = 3838.5386
$ws.Range("J136").Value = 8953.333000000001
$ws.Range("K136").Value = 11515.6158
$ws.Range("L136").Value = 26859.999
$ws.Range("M136").Value = -8965.6158
$ws.Range("N136").Value = -31959.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1047.8572
$ws.Range("I100").Value = 283
$ws.Range("K100").Value = 566
$ws.Range("M100").Value = -25

$ws.Range("H113").Value = 967.6
$ws.Range("I113").Value = 1142.6
$ws.Range("J113").Value = 442.6
$ws.Range("K113").Value = 3427.8
$ws.Range("L113").Value = 1327.8
$ws.Range("M113").Value = -1257.8
$ws.Range("N113").Value = -5667.8

$ws.Range("H126").Value = 2213.2
$ws.Range("I126").Value = 1453.8572
$ws.Range("K126").Value = 4361.571599999999
$ws.Range("M126").Value = -1891.571599999999

$ws.Range("H132").Value = 2005.8649
$ws.Range("J132").Value = 2813.625
$ws.Range("L132").Value = 8440.875
$ws.Range("N132").Value = -13500.875
